$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Weekly crime table updates (rows 14-29) ---

# Row 14
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 2

# Row 15
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("N15").Value = -86.666666666666

# Row 16
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -87.5
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -46.666666666666
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = -62.857142857142
$ws.Range("M16").Value = -59.375
$ws.Range("N16").Value = -91.925465838509

# Row 17
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 48
$ws.Range("I17").Value = 57
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 14
$ws.Range("L17").Value = 21.27659574468
$ws.Range("M17").Value = 83.870967741935
$ws.Range("N17").Value = -37.362637362637

# Row 18
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 3
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 9
$ws.Range("K18").Value = -47.058823529411
$ws.Range("L18").Value = 28.571428571428
$ws.Range("M18").Value = -65.384615384615
$ws.Range("N18").Value = -95.477386934673

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -28.571428571428
$ws.Range("I19").Value = 42
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = -28.813559322033
$ws.Range("L19").Value = -28.813559322033
$ws.Range("M19").Value = 2.439024390243
$ws.Range("N19").Value = -41.666666666666

# Row 20
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -46.666666666666
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = -26.315789473684
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = -44
$ws.Range("N20").Value = -93.26923076923

# Row 21
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = -15.306122448979
$ws.Range("I21").Value = 137
$ws.Range("J21").Value = 183
$ws.Range("K21").Value = -25.136612021857
$ws.Range("L21").Value = 0.735294117647
$ws.Range("M21").Value = -14.375
$ws.Range("N21").Value = -81.684491978609

# Row 23
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -44.444444444444
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = -55.555555555555
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 300

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 68.421052631578
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 25.925925925925
$ws.Range("I24").Value = 178
$ws.Range("J24").Value = 161
$ws.Range("K24").Value = 10.55900621118
$ws.Range("L24").Value = 64.814814814814
$ws.Range("M24").Value = 0.564971751412

# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 15.384615384615
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -9.090909090909
$ws.Range("I25").Value = 83
$ws.Range("J25").Value = 95
$ws.Range("K25").Value = -12.631578947368
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -39.855072463768

# Row 26
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 6
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "0"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "***.*"
$ws.Range("I26").Value = 8
$ws.Range("K26").Value = 700
$ws.Range("L26").Value = 166.666666666667

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("I27").Value = 15
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 150

# Row 28
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("G28").Value = 1
$ws.Range("H28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H28").Value = 100
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("J28").Value = 1
$ws.Range("K28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("K28").Value = 200
$ws.Range("N28").Value = -78.571428571428

# Row 29
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E29").Value = -100
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("G29").Value = 1
$ws.Range("H29").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H29").Value = 100
$ws.Range("J29").NumberFormat = "#,##0"
$ws.Range("J29").Value = 1
$ws.Range("K29").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("K29").Value = 200
$ws.Range("N29").Value = -75
